$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the ZGN (Zhongshan, China) row. It is currently row 242.
# Deleting the entire row shifts all subsequent rows up by one,
# which matches the diff (dimension shrinks from H332 to H331).
$ws.Rows.Item(242).Delete()
